$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new price values in column E for rows 3 and 4
$ws.Range("E3").Value = 6
$ws.Range("E4").Value = 6.25

# Move the active selection from D5 to E5
$ws.Range("E5").Select()
